$d = $word.ActiveDocument

# 1) "Test " + "de integracion" were two separate runs with identical
#    run-formatting; collapse them into a single run "Test de integracion".
#    Re-applying Find/Replace with the same text (case-sensitive, so the
#    lookalike text elsewhere in the document is left untouched) makes Word
#    rewrite the matched range as one run using the first run's formatting,
#    which is exactly the merge the diff shows.
$d.Content.Find.Execute(
    "Test de integracion", $true, $false, $false, $false, $false,
    $true, 1, $false, "Test de integracion", 2) | Out-Null

# 2) Likewise "Test " + "Suite" (two runs) -> "Test Suite" (one run).
#    MatchCase=$true keeps this from touching the unrelated lowercase
#    "JUnit test suite ..." sentence a little further down.
$d.Content.Find.Execute(
    "Test Suite", $true, $false, $false, $false, $false,
    $true, 1, $false, "Test Suite", 2) | Out-Null

# 3) Append two new blank paragraphs (centered, bold, underlined, Arial 36)
#    right after the last paragraph ("Patron Cadena de Responsabilidad"),
#    before the section properties. InsertXML lets us add the paragraph
#    with only a <w:pPr> (carrying the mark's run-formatting in rPr) and no
#    run at all, matching the target markup exactly.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/><w:u w:val="single"/><w:lang w:eastAsia="es-AR"/></w:rPr></w:pPr></w:p>'

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($newParaXml) | Out-Null

$endRange2 = $d.Range($d.Content.End, $d.Content.End)
$endRange2.InsertXML($newParaXml) | Out-Null
